$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add date value to C3 (date column), formatted like C2/C4 (mm-dd-yy date format)
$ws.Range("C4").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("C3").Value = 41256

# Add value to D4 (zs_dz column), row 4
$ws.Range("D4").Value = 3

# Update selection to E4
$ws.Range("E4").Select()
